# Auto-generated edit script for cryptos.xlsx update
# Updates Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.543.13"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.873.40"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.31"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4757"
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2908"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06489"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.84"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07749"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7370"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "1.873.02"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.86"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.171"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.95"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "30.610.43"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.18"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007482"
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "2.121.28"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.201"
$ws.Range("E23").Value = "  -1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.164"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.171"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.52"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.75"
$ws.Range("E27").Value = "  -2.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.903"
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09862"
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.340"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.500"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.257"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.083"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04782"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6925"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01849"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.753"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.47"
$ws.Range("E41").Value = "  +4.15%  "
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4170"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8348"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.43"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.351"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.22"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.954"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "915.51"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05668"
